$d = $word.ActiveDocument

# The "Highlights" section used to contain a single paragraph listing the
# four highlight color codes concatenated together. Split it into four
# separate paragraphs, each annotated with its count.
$d.Content.Find.Execute(
    "#7cc867#fb5b89#c885da#f9cd59",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "#7cc867: 48^p#fb5b89: 43^p#c885da: 19^p#f9cd59: 55",
    2
)
